$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DOMA-7423: updated meter import example - diversify the "Unit type" example
# values for the first address (rows 3-6) so the sample file demonstrates
# the full range of supported unit types.
$ws.Range("C3").Value = "Parking place"
$ws.Range("C4").Value = "Apartment"
$ws.Range("C5").Value = "Warehouse unit"
$ws.Range("C6").Value = "Commercial unit"
